$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.797.10'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.88%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.814.83'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.16%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '277.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -8.18%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5096'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -5.03%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3520'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.39%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.62'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.88%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06667'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -7.26%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.01'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -7.33%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.8293'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -6.62%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07859'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.75%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.823.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.65%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -6.58%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -4.25%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008028'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -6.11%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '25.844.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.88%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.723'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.22%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -6.31%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.03%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Monero'
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.77%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.193'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.82%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.672'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.69%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.06'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -5.36%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '109.38'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.07%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.351'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.95%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -8.05%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.98%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04902'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.52%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7290'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -9.47%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.05%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.877'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.0000'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.138'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.02%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.98%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5198'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -13.70%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.26%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9569'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -10.59%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.50%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '111.55'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.05%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.011'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -9.53%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4557'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -11.44%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -8.71%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.49%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.226'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -7.09%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -8.00%  '
$ws.Range('E51').Style = 'Normal'
